$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The yearly report is shifting its rolling 5-year window forward by one year:
# the oldest year (1396/12) is dropped, the four remaining years shift one
# column to the left (E<-F, F<-G, G<-H, H<-I), and the newest year (1401/12)
# is appended in column I together with its figures.

# Header rows that hold the "twelve months ended <date>" column captions.
$headerRows = @(8, 24)
foreach ($r in $headerRows) {
    $fVal = $ws.Range("F$r").Value2
    $gVal = $ws.Range("G$r").Value2
    $hVal = $ws.Range("H$r").Value2
    $iVal = $ws.Range("I$r").Value2

    $ws.Range("E$r").Value2 = $fVal
    $ws.Range("F$r").Value2 = $gVal
    $ws.Range("G$r").Value2 = $hVal
    $ws.Range("H$r").Value2 = $iVal
    $ws.Range("I$r").Value2 = "دوازده ماهه منتهی به 1401/12"
}

# Data rows: shift the four most recent of the old figures left by one column
# and place the freshly reported 1401/12 figure in column I.
$rowsNewValues = @{
    10 = 6428095   # هزینه حمل و نقل و انتقال
    11 = 0         # هزینه خدمات پس از فروش
    12 = 66820     # حق العمل و کمیسیون فروش
    13 = 0         # هزینه تبلیغات
    14 = 41430     # هزینه مواد مصرفی
    15 = 0         # هزینه انرژی (آب، برق، گاز و سوخت)
    16 = 13324     # هزینه استهلاک
    17 = 788030    # هزینه حقوق و دستمزد
    18 = 0         # هزینه مطالبات مشکوک الوصول
    19 = 2091584   # سایر هزینه ها
    20 = 9429283   # جمع
    26 = 220       # تعداد پرسنل غیر تولیدی شرکت
    27 = 85        # تعداد پرسنل تولیدی شرکت
}

foreach ($r in $rowsNewValues.Keys) {
    $fVal = $ws.Range("F$r").Value2
    $gVal = $ws.Range("G$r").Value2
    $hVal = $ws.Range("H$r").Value2
    $iVal = $ws.Range("I$r").Value2

    $ws.Range("E$r").Value2 = $fVal
    $ws.Range("F$r").Value2 = $gVal
    $ws.Range("G$r").Value2 = $hVal
    $ws.Range("H$r").Value2 = $iVal
    $ws.Range("I$r").Value2 = $rowsNewValues[$r]
}

Write-Host "Done applying updates"
